$d = $word.ActiveDocument

# Update the date heading paragraph.
$d.Content.Find.Execute("2025-11-03 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-11-04 Tuesday", 2)

# Update the division problems in the table by cell coordinates, so the
# duplicate "89" / "88" values in the grid don't get mixed up.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="88÷7="},
    @{Row=1;  Col=2; Text="69÷8="},
    @{Row=1;  Col=3; Text="81÷5="},
    @{Row=1;  Col=4; Text="13÷8="},
    @{Row=1;  Col=5; Text="35÷6="},

    @{Row=5;  Col=1; Text="27÷3="},
    @{Row=5;  Col=2; Text="15÷9="},
    @{Row=5;  Col=3; Text="28÷6="},
    @{Row=5;  Col=4; Text="99÷2="},
    @{Row=5;  Col=5; Text="17÷7="},

    @{Row=9;  Col=1; Text="92÷4="},
    @{Row=9;  Col=2; Text="19÷6="},
    @{Row=9;  Col=3; Text="52÷3="},
    @{Row=9;  Col=4; Text="30÷2="},
    @{Row=9;  Col=5; Text="98÷9="},

    @{Row=13; Col=1; Text="39÷2="},
    @{Row=13; Col=2; Text="94÷8="},
    @{Row=13; Col=3; Text="81÷9="},
    @{Row=13; Col=4; Text="19÷4="},
    @{Row=13; Col=5; Text="52÷9="},

    @{Row=17; Col=1; Text="79÷7="},
    @{Row=17; Col=2; Text="11÷4="},
    @{Row=17; Col=3; Text="98÷7="},
    @{Row=17; Col=4; Text="30÷6="},
    @{Row=17; Col=5; Text="21÷8="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
